$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 3429.0908
$ws.Range("I80").Value = 2502.3333
$ws.Range("J80").Value = 4541.2
$ws.Range("K80").Value = 7506.999899999999
$ws.Range("L80").Value = 13623.6
$ws.Range("M80").Value = -6508.999899999999
$ws.Range("N80").Value = -15619.6
$ws.Range("H81").Value = 29412.75
$ws.Range("J81").Value = 29412.75
$ws.Range("L81").Value = 29412.75
$ws.Range("N81").Value = -31408.75
$ws.Range("H83").Value = 3429.0908
$ws.Range("I83").Value = 2502.3333
$ws.Range("J83").Value = 4541.2
$ws.Range("K83").Value = 22520.9997
$ws.Range("L83").Value = 40870.8
$ws.Range("M83").Value = -17528.9997
$ws.Range("N83").Value = -50854.8
$ws.Range("H84").Value = 29412.75
$ws.Range("J84").Value = 29412.75
$ws.Range("L84").Value = 88238.25
$ws.Range("N84").Value = -98222.25
$ws.Range("H88").Value = 3372.8948
$ws.Range("I88").Value = 3967.1667
$ws.Range("J88").Value = 3098.6155
$ws.Range("K88").Value = 3967.1667
$ws.Range("L88").Value = 3098.6155
$ws.Range("M88").Value = -3561.1667
$ws.Range("N88").Value = -3910.6155
$ws.Range("H91").Value = 3372.8948
$ws.Range("I91").Value = 3967.1667
$ws.Range("J91").Value = 3098.6155
$ws.Range("K91").Value = 3967.1667
$ws.Range("L91").Value = 3098.6155
$ws.Range("M91").Value = -2563.1667
$ws.Range("N91").Value = -5906.6155
$ws.Range("H106").Value = 8700.799999999999
$ws.Range("I106").Value = 3996.6667
$ws.Range("J106").Value = 10716.857
$ws.Range("K106").Value = 3996.6667
$ws.Range("L106").Value = 10716.857
$ws.Range("M106").Value = -3365.6667
$ws.Range("N106").Value = -11978.857
$ws.Range("H138").Value = 2819539.5
$ws.Range("I138").Value = 1640.5807
$ws.Range("J138").Value = 5003411
$ws.Range("K138").Value = 4921.742099999999
$ws.Range("L138").Value = 15010233
$ws.Range("M138").Value = 218.2579000000005
$ws.Range("N138").Value = -15020513
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 24990.977
$ws.Range("I74").Value = 38053.332
$ws.Range("J74").Value = 1478.7333
$ws.Range("K74").Value = 38053.332
$ws.Range("L74").Value = 1478.7333
$ws.Range("M74").Value = -37179.332
$ws.Range("N74").Value = -3226.7333
$ws.Range("H77").Value = 24990.977
$ws.Range("I77").Value = 38053.332
$ws.Range("J77").Value = 1478.7333
$ws.Range("K77").Value = 190266.66
$ws.Range("L77").Value = 7393.6665
$ws.Range("M77").Value = -185898.66
$ws.Range("N77").Value = -16129.6665
$ws.Range("H88").Value = 2210.1765
$ws.Range("I88").Value = 2020.75
$ws.Range("J88").Value = 2378.5557
$ws.Range("K88").Value = 2020.75
$ws.Range("L88").Value = 2378.5557
$ws.Range("M88").Value = -1614.75
$ws.Range("N88").Value = -3190.5557
$ws.Range("H91").Value = 2210.1765
$ws.Range("I91").Value = 2020.75
$ws.Range("J91").Value = 2378.5557
$ws.Range("K91").Value = 2020.75
$ws.Range("L91").Value = 2378.5557
$ws.Range("M91").Value = -616.75
$ws.Range("N91").Value = -5186.5557
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2287.394
$ws.Range("J62").Value = 2354
$ws.Range("L62").Value = 2354
$ws.Range("N62").Value = -3602
$ws.Range("H65").Value = 2287.394
$ws.Range("J65").Value = 2354
$ws.Range("L65").Value = 11770
$ws.Range("N65").Value = -18010
$ws.Range("H105").Value = 949.4
$ws.Range("I105").Value = 947.2727
$ws.Range("J105").Value = 955.25
$ws.Range("K105").Value = 947.2727
$ws.Range("L105").Value = 955.25
$ws.Range("M105").Value = 799.7273
$ws.Range("N105").Value = -4449.25
$ws.Range("H132").Value = 1060369.6
$ws.Range("I132").Value = 2135.68
$ws.Range("J132").Value = 3705954.5
$ws.Range("K132").Value = 6407.039999999999
$ws.Range("L132").Value = 11117863.5
$ws.Range("M132").Value = -3877.039999999999
$ws.Range("N132").Value = -11122923.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 914.48486
$ws.Range("J131").Value = 916.14734
$ws.Range("L131").Value = 2748.44202
$ws.Range("N131").Value = -12828.44202
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2449.8
$ws.Range("I80").Value = 2361.389
$ws.Range("K80").Value = 2361.389
$ws.Range("M80").Value = -1363.389
$ws.Range("H83").Value = 2449.8
$ws.Range("I83").Value = 2361.389
$ws.Range("K83").Value = 11806.945
$ws.Range("M83").Value = -6814.945
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 13316.667
$ws.Range("I29").Value = 8900
$ws.Range("J29").Value = 14200
$ws.Range("K29").Value = 8900
$ws.Range("L29").Value = 14200
$ws.Range("M29").Value = -8605
$ws.Range("N29").Value = -14790
$ws.Range("H68").Value = 7328
$ws.Range("I68").Value = 18731.666
$ws.Range("J68").Value = 3051.625
$ws.Range("K68").Value = 18731.666
$ws.Range("L68").Value = 3051.625
$ws.Range("M68").Value = -17982.666
$ws.Range("N68").Value = -4549.625
$ws.Range("H71").Value = 7328
$ws.Range("I71").Value = 18731.666
$ws.Range("J71").Value = 3051.625
$ws.Range("K71").Value = 93658.33
$ws.Range("L71").Value = 15258.125
$ws.Range("M71").Value = -89914.33
$ws.Range("N71").Value = -22746.125
$ws.Range("H136").Value = 1577.122
$ws.Range("I136").Value = 1056.75
$ws.Range("J136").Value = 2311.7646
$ws.Range("K136").Value = 3170.25
$ws.Range("L136").Value = 6935.293799999999
$ws.Range("M136").Value = -620.25
$ws.Range("N136").Value = -12035.2938
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1125.25
$ws.Range("I81").Value = 1163.8182
$ws.Range("J81").Value = 1040.4
$ws.Range("K81").Value = 2327.6364
$ws.Range("L81").Value = 2080.8
$ws.Range("M81").Value = -1266.6364
$ws.Range("N81").Value = -4202.8
$ws.Range("H84").Value = 1125.25
$ws.Range("I84").Value = 1163.8182
$ws.Range("J84").Value = 1040.4
$ws.Range("K84").Value = 11638.182
$ws.Range("L84").Value = 10404
$ws.Range("M84").Value = -6334.181999999999
$ws.Range("N84").Value = -21012
$ws.Range("H122").Value = 107530.62
$ws.Range("I122").Value = 26869
$ws.Range("J122").Value = 201635.83
$ws.Range("K122").Value = 80607
$ws.Range("L122").Value = 604907.49
$ws.Range("M122").Value = -78157
$ws.Range("N122").Value = -609807.49
$ws.Range("H126").Value = 1050
$ws.Range("I126").Value = 926.9231
$ws.Range("J126").Value = 1850
$ws.Range("K126").Value = 2780.7693
$ws.Range("L126").Value = 5550
$ws.Range("M126").Value = -310.7692999999999
$ws.Range("N126").Value = -10490
